$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footer "first page" (file word/footer1.xml) ---------------------------
# Pearson logo: wp:docPr/pic:cNvPr id="3"/"0", name "image1.png" -> "image2.png"
$fFirst = $sec.Footers(2)
$ishp = $fFirst.Range.InlineShapes(1)
$shp = $ishp.ConvertToShape()
$shp.Name = "image2.png"
$shp.ConvertToInlineShape() | Out-Null

# --- Footer "default" (file word/footer2.xml) -------------------------------
# Pearson logo: wp:docPr/pic:cNvPr id="2"/"0", name "image1.png" -> "image2.png"
$fDefault = $sec.Footers(1)
$ishp = $fDefault.Range.InlineShapes(1)
$shp = $ishp.ConvertToShape()
$shp.Name = "image2.png"
$shp.ConvertToInlineShape() | Out-Null

# --- Header "first page" (file word/header1.xml) ----------------------------
# BTec logo: wp:docPr/pic:cNvPr id="1"/"0", name "image2.jpg" -> "image1.jpg"
$hFirst = $sec.Headers(2)
$ishp = $hFirst.Range.InlineShapes(1)
$shp = $ishp.ConvertToShape()
$shp.Name = "image1.jpg"
$shp.ConvertToInlineShape() | Out-Null
